$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newRow = 91

# Leading apostrophe forces the date-looking text to stay as text instead of
# being auto-converted into a real date serial number by Excel.
$ws.Cells.Item($newRow, 1).Value = "'2026/02/09"
$ws.Cells.Item($newRow, 2).Value = "逃离鸭科夫"
$ws.Cells.Item($newRow, 3).Value = 1185

# Reset to the "Normal" cell style so the quote-prefix entry doesn't leave a
# stray formatting flag behind, then re-apply the centered alignment used by
# every other data row in the sheet.
$ws.Cells.Item($newRow, 1).Style = "Normal"
$ws.Cells.Item($newRow, 2).Style = "Normal"
$ws.Cells.Item($newRow, 3).Style = "Normal"

$srcRange = $ws.Range("A90:C90")
$dstRange = $ws.Range("A91:C91")
$dstRange.HorizontalAlignment = $srcRange.HorizontalAlignment
$dstRange.VerticalAlignment = $srcRange.VerticalAlignment
